## B6-PowerPoint.pptx edit
##  1. Re-apply a table style (by GUID) to the three data tables in the deck.
##  2. Swap the theme applied to the slide master with the theme applied to
##     the notes master (Integral <-> Office Theme), by swapping every
##     theme colour between the two embedded themes.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Tables: switch from the custom "Table_0" style to the built-in table
#    style {1D2D30AA-F715-4C8B-8720-6D4A47C31290}.
# ---------------------------------------------------------------------------
$newTableStyleId = "{1D2D30AA-F715-4C8B-8720-6D4A47C31290}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Themes: the slide master currently carries the "Integral" theme while
#    the notes master carries the "Office Theme" theme. Swap the two themes'
#    colour schemes (the only parts that differ between them) so the slide
#    master ends up with the Office colours and the notes master ends up
#    with the Integral colours.
# ---------------------------------------------------------------------------

# Helper: convert a 0xRRGGBB value into the BGR long used by ColorFormat.RGB.
function ToBgr([int]$rrggbb) {
    $r = ($rrggbb -shr 16) -band 0xFF
    $g = ($rrggbb -shr 8) -band 0xFF
    $b = $rrggbb -band 0xFF
    return ($b * 65536) + ($g * 256) + $r
}

# Colour order matches ThemeColorScheme.Item(1..12):
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    0x000000, 0xFFFFFF, 0x44546A, 0xE7E6E6,
    0x5B9BD5, 0xED7D31, 0xA5A5A5, 0xFFC000, 0x4472C4, 0x70AD47,
    0x0563C1, 0x954F72
)
$integralThemeColors = @(
    0x000000, 0xFFFFFF, 0x454551, 0xD8D9DC,
    0xE32D91, 0xC830CC, 0x4EA6DC, 0x4775E7, 0x8971E1, 0xD54773,
    0x6B9F25, 0x8C8C8C
)

$slideMasterColors = $p.SlideMaster.Theme.ThemeColorScheme
$notesMasterColors = $p.NotesMaster.Theme.ThemeColorScheme

for ($k = 1; $k -le 12; $k++) {
    $slideMasterColors.Item($k).RGB = ToBgr($officeThemeColors[$k - 1])
    $notesMasterColors.Item($k).RGB = ToBgr($integralThemeColors[$k - 1])
}
